$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Text change: "Ready for handoff" -> "In Translation" -------------------
# Update every cell that currently shows the old status so the shared string
# table collapses back onto a single (renamed) entry instead of leaving an
# orphaned "Ready for handoff" string behind.
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column width change -----------------------------------------------------
# Overview!E:F and zh-cn!C / de-de!C narrow from ~17.22 chars to ~13.41 chars.
$newWidth = 12.5
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
